# Update gh-pages output data (合肥-漫展信息.xlsx)
$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibition): refresh "想去人数" (interested-count) figures ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 555
$ws1.Range("F4").Value = 49
$ws1.Range("F9").Value = 436
$ws1.Range("F10").Value = 3509

# --- Sheet "演出" (Performance): the 2024-08-25 CrossingX event has expired, remove it ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Rows.Item(2).Delete()

# --- Sheet "全部类型" (All types): same expired CrossingX event, remove it ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Rows.Item(3).Delete()
